# Apply the crypto price-ticker refresh described in the commit diff.
# Column D ("Price") holds numeric-looking text (e.g. "243.74", "30.330.59")
# that must stay text, matching the workbook's inline-string cells. A leading
# apostrophe forces text entry (like typing it in Excel), then resetting the
# cell Style back to "Normal" strips the quote-prefix style Excel applies,
# so the cell ends up plain text with the default style - exactly like the original.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.330.59"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Value = "'1.872.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'243.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.59%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D8").Value = "'0.2877"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.82%  "

$ws.Range("D9").Value = "'0.06448"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").Value = "'22.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("D11").Value = "'0.07774"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("D12").Value = "'1.874.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.05%  "

$ws.Range("D13").Value = "'95.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").Value = "'0.7212"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "

$ws.Range("D15").Value = "'5.131"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.78%  "

$ws.Range("D16").Value = "'279.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("D17").Value = "'30.320.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.82%  "

$ws.Range("D18").Value = "'13.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "'0.000007441"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.37%  "

$ws.Range("D21").Value = "'2.118.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("D22").Value = "'0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").Value = "'5.238"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.55%  "

$ws.Range("D24").Value = "'6.238"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").Value = "'163.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.25%  "

$ws.Range("D26").Value = "'9.052"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.48%  "

$ws.Range("D27").Value = "'18.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.35%  "

$ws.Range("D28").Value = "'1.880"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.19%  "

$ws.Range("D29").Value = "'1.316"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.16%  "

$ws.Range("D30").Value = "'0.09574"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.15%  "

$ws.Range("E32").Value = "  -0.75%  "

$ws.Range("E33").Value = "  +0.09%  "

$ws.Range("D34").Value = "'0.04804"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.65%  "

$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").Value = "'0.6885"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.75%  "

$ws.Range("D38").Value = "'0.01875"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.46%  "

$ws.Range("E39").Value = "  +1.96%  "

$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("D41").Value = "'74.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.37%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.942"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.71%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4235"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.69%  "

$ws.Range("D44").Value = "'0.9991"
$ws.Range("D44").Style = "Normal"

$ws.Range("E45").Value = "  -1.31%  "

$ws.Range("D46").Value = "'100.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.84%  "

$ws.Range("D47").Value = "'9.580"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.65%  "

$ws.Range("D48").Value = "'35.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.57%  "

$ws.Range("D49").Value = "'6.921"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.54%  "

$ws.Range("D50").Value = "'899.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "

$ws.Range("D51").Value = "'0.05718"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.91%  "
